# Actualización 11 de Mayo - Tarde
# Add a new rescued-student row to the "Rescatables" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 19330051920227
$ws.Range("B2").Value = "CHORA"
$ws.Range("C2").Value = "LOPEZ"
$ws.Range("D2").Value = "GABRIEL ALEJANDRO"
$ws.Range("E2").Value = "DESARROLLA APLICACIONES QUE SE EJECUTAN EN EL CLIENTE"
$ws.Range("F2").Value = "4APM"
$ws.Range("G2").Value = 2
